$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 125, shifting existing rows 125-133 down to 126-134
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the latest weekly data point
$ws.Range("A125").Value = 11
$ws.Range("B125").Value = "Vega Monumental Concepción"
$ws.Range("C125").Value = "Bíobío"
$ws.Range("D125").Value = 44931
$ws.Range("E125").Value = 8
$ws.Range("F125").Value = 100112001
$ws.Range("G125").Value = "Berenjena"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 15000
$ws.Range("L125").Value = 16000
$ws.Range("M125").Value = 15500
$ws.Range("N125").Value = "`$/caja 60 unidades"
$ws.Range("O125").Value = "Región de Arica y Parinacota"
$ws.Range("P125").Value = 258
$ws.Range("Q125").Value = 60
$ws.Range("R125").Value = "Hortaliza"
